# Applies the 9-May-2023 cryptos-list refresh: updates the Price (D) and
# Volume(1h) (E) columns for the rows whose quoted figures moved.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.664.04"
$ws.Range("E2").Value = "  -1.30%  "
$ws.Range("D3").Value = "1.851.43"
$ws.Range("E3").Value = "  -1.03%  "
$ws.Range("E4").Value = "  -0.37%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.56%  "
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4241"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.00%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3644"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.51"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07299"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8763"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.75"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.18%  "
$ws.Range("D13").Value = "1.840.58"
$ws.Range("E13").Value = "  -7.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.346"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.42%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.531"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.98%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06889"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.28%  "
$ws.Range("E17").Value = "  -0.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "79.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008909"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("E21").Value = "  -2.21%  "
$ws.Range("D22").Value = "27.688.29"
$ws.Range("E22").Value = "  -1.21%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.991"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.27%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.20%  "
$ws.Range("D25").Value = "2.067.73"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.983"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "152.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.96%  "
$ws.Range("E28").Value = "  +2.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "121.91"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.278"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.44%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.880"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.52%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08872"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.30%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7681"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.71%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.568"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.85%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.973"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.61%  "
$ws.Range("E36").Value = "  -5.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9996"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.56%  "
$ws.Range("E38").Value = "  -1.88%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01936"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.68%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.811"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.52%  "
$ws.Range("E42").Value = "  -2.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.889"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1650"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.311"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06546"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.47%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4766"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.99%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.38"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.60"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.9997"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("E51").Value = "  -2.44%  "
